$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New/updated long-form text blocks (content of this commit)
# ---------------------------------------------------------------------------

$objetivos = @"
1. Introduzir e discutir conceitos e técnicas estatísticas para controle e melhoria da qualidade de produtos fabricados e processos de fabricação;2 Fornecer subsídios para que o aluno tenha condições de utilizar essas técnicas e conceitos na sua vida profissional futura.
"@

$docentes = @"
5840917 - Fabricio Maciel Gomes
"@

$resumo = @"
Sistemas de Medição, Fundamentos do Controle Estatístico da Qualidade e do Processo, Gráficos de Controle por Variáveis, Gráficos de Controle por Atributos, Gráficos de Controle para Processos Auto-correlacionados, Analise de Capacidade do Processo, Inspeção da Qualidade, Estudos de casos.
"@

$programa = @"
1. Sistemas de Medição.1.1. Planejamento do Sistema de Medição;1.2. Impacto da Variabilidade do Sistema de Medição no Produto;1.3. Sistemas de Medição por Atributos;1.4. Tendência e Linearidade;1.5. Análise de Repetitividade e Reprodutibilidade;2. Fundamentos do Controle Estatístico da Qualidade e do Processo.2.1. Importância do Controle Estatístico da Qualidade e do Processo;2.2. Naturezas das Variações;2.3. Causas Comuns e Causas Especiais de Variações;3. Gráficos de Controle por Variáveis3.1. Gráficos de Controle por Médias;3.2. Gráficos de Controle por Amplitude;3.3. Gráficos de Controle por Desvio Padrão;3.4. Análise de Desempenho dos Gráficos de Controle por Variáveis;4. Gráficos de Controle por Atributos4.1. Gráficos de Controle por Número de Não Conformidades;4.2. Gráficos de Controle por Fração Não Conforme;4.3. Gráficos de Controle por Número de Defeitos4.4. Gráficos de Controle por Não Conformidades por Amostra;5. Gráficos de Controle para Processos Auto-correlacionados5.1. Gráficos de Controle por Amplitude Móvel;5.2. Gráficos de Controle por Soma Acumulada (CUSUM).5.3. Gráficos de Controle por Média Móvel Ponderada Exponencialmente (EWMA)6. Analise de Capacidade do Processo6.1. Índices de Capacidade do Processo;6.2. Índices de Performance do Processo; 7. Inspeção da Qualidade7.1. Planos de Amostragem7.2. Inspeção para Aceitação;7.3. Inspeção Retificadora;8. Estudos de casos
"@

$metodo = @"
Aulas expositivas teóricas, aulas práticas, aulas de laboratório, aulas de exercícios.
"@

$criterio = @"
MF = (0,5*P1 + 0,5*P2), onde P1 e P2 são provas.
"@

$norma = @"
Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.
"@

$biblio = @"
1. COSTA, A.F.B., EPPRECHT, E.K., CARPINETTI, L.C.R., Controle Estatístico da Qualidade, 2ª ed., Editora Atlas, 2005.
2. MONTGOMERY, D.C., Introdução ao Controle Estatístico da Qualidade, 4ª ed., Livros Técnicos e Científicos, 2004.
3. GRANT, E., LEAVENWORTH, R., Statistical Quality Control, 7ªed., McGraw-Hill, 1996.
4. WERKENA, M.C.C., Ferramentas Estatísticas Básicas para o Gerenciamento de Processos, Editora FCO, 1996.
"@

# ---------------------------------------------------------------------------
# 1. Row 10 ("Objetivos:") previously leaked the docente string into B/C;
#    fix it to the real course-objectives paragraph.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------------
# 2. Insert a brand new row 13 (pushes everything below it down by one) that
#    carries the "Docentes responsáveis" value, which used to (incorrectly)
#    sit on the "Método:" row further down.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).Clear()

$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)

$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes

# ---------------------------------------------------------------------------
# 3. Former row 13 ("Programa resumido:" / "Semestral") is now row 14; give
#    it the real short-syllabus paragraph.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = $resumo
$ws.Range("C14").Value = $resumo

# ---------------------------------------------------------------------------
# 4. Former row 15 ("Programa:" / "01/01/2018") is now row 16; give it the
#    real full syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# ---------------------------------------------------------------------------
# 5. Former row 18 ("Método:", which had the docente string leaked into it)
#    is now row 19; fix it to the real teaching-method text.
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# ---------------------------------------------------------------------------
# 6. Former row 19 ("Critério:", which had the método text leaked into it)
#    is now row 20; fix it to the real grading-criteria formula.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# ---------------------------------------------------------------------------
# 7. Former row 20 ("Norma de recuperação:", which had the critério text
#    leaked into it) is now row 21; fix it to the real recovery-grade text.
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# ---------------------------------------------------------------------------
# 8. Former row 21 ("Bibliografia:", which had the norma text leaked into
#    it) is now row 22; fix it to the real bibliography text.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
